$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for the "loading_percent" results table, columns B:F, rows 2:25
# (re-computed results after adding Rac, Dext, kg to the underlying model)
$values = @(
    @(26.4234750373145,19.99886229031786,9.388288595591627,30.77377411776597,63.56182861480239),
    @(24.38484812770647,18.46397124672549,8.849519637809431,28.31226926357537,59.19036081192131),
    @(23.08887446149555,17.49128058524141,8.517092694060617,26.76069239565595,56.45050333521705),
    @(22.5488165671451,17.08658420349867,8.381070846560407,26.11693298543749,55.31894217127872),
    @(22.45839791254524,17.01886368656673,8.358447793206549,26.00930803074004,55.13011678721339),
    @(23.08164037263802,17.48585724081155,8.515260655957015,26.75205855939405,56.43530467096087),
    @(25.72920989287836,19.47544971045623,9.202748790060806,29.93249018775551,62.06535951554753),
    @(30.6151209524473,23.17640264765306,10.54862083012398,35.92794154082256,72.75185367649347),
    @(34.10049648182896,25.84552666713169,11.55880171813866,40.33092135322983,80.76514898864107),
    @(35.69037928621933,27.07283950664322,12.03142000939431,42.38294602322394,84.44266199602254),
    @(36.29636231448981,27.54243372989783,12.21335939014254,43.17331847737465,85.84726717030948),
    @(36.16561187650134,27.44102345508451,12.17402115518748,43.00237729302152,85.54407671466608),
    @(35.74013170463823,27.11135533341069,12.04632091510747,42.44765844908599,84.55792630318818),
    @(35.48014364190327,26.91016096761659,11.96852738362011,42.1098404086239,83.95571276964839),
    @(33.99693745489056,25.76581099699495,11.52826172298007,40.19828504984866,80.52603259143342),
    @(33.09015736458159,25.06893172664895,11.26217231165911,39.04192034363454,78.43472975276252),
    @(32.56876172698243,24.66911465548086,11.11027157901295,38.38094843792835,77.23433467117515),
    @(32.39217126398298,24.53384343894225,11.05901096618923,38.15771199574974,76.82814258416828),
    @(33.18664885298153,25.14299322364008,11.29037273589011,39.16455127020188,78.65705230527773),
    @(35.86496630445502,27.20802630276621,12.08373806170879,42.61016862356581,84.84718327092956),
    @(37.64089250772431,28.58817570060622,12.62043690929643,44.94478676273031,88.96865687882396),
    @(36.68925693826737,27.84746056301008,12.33183565613295,43.68836261057393,86.75872425500344),
    @(33.14302551706285,25.10950764703843,11.27762004292854,39.10909830811917,78.55653480363642),
    @(29.31887174341284,22.19033144700976,10.18291235798907,34.31926875976241,69.87664922851232)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($row, 2 + $j).Value = $rowVals[$j]
    }
}
